$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -1
$ws.Range("F3").Value = -3
$ws.Range("F4").Value = -8
$ws.Range("F5").Value = -1
